# Re-doing global M2 module
# Applies updated M2_Len/FX_Len counts and refreshed M2/FX first/last date serials
# for the Long28_DataComp sheet (Global_M2 liquidity datasums).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("C2").Value = 360
    $ws.Range("F2").Value = 45992
    $ws.Range("G2").Value = 30865
    $ws.Range("H2").Value = 46055
    # Row 3
    $ws.Range("E3").Value = 30803
    $ws.Range("F3").Value = 45992
    # Row 4
    $ws.Range("E4").Value = 30803
    $ws.Range("F4").Value = 45992
    $ws.Range("G4").Value = 30865
    $ws.Range("H4").Value = 46055
    # Row 5
    $ws.Range("E5").Value = 30803
    $ws.Range("F5").Value = 45992
    $ws.Range("G5").Value = 30865
    $ws.Range("H5").Value = 46055
    # Row 6
    $ws.Range("G6").Value = 30865
    $ws.Range("H6").Value = 46055
    # Row 7
    $ws.Range("E7").Value = 30773
    $ws.Range("F7").Value = 45962
    $ws.Range("G7").Value = 30865
    $ws.Range("H7").Value = 46055
    # Row 8
    $ws.Range("D8").Value = 436
    $ws.Range("H8").Value = 46055
    # Row 9
    $ws.Range("G9").Value = 30865
    $ws.Range("H9").Value = 46055
    # Row 10
    $ws.Range("E10").Value = 30803
    $ws.Range("F10").Value = 45992
    $ws.Range("G10").Value = 30865
    $ws.Range("H10").Value = 46055
    # Row 11
    $ws.Range("E11").Value = 30773
    $ws.Range("F11").Value = 45962
    $ws.Range("G11").Value = 30865
    $ws.Range("H11").Value = 46055
    # Row 12
    $ws.Range("C12").Value = 397
    $ws.Range("D12").Value = 378
    $ws.Range("F12").Value = 45992
    $ws.Range("H12").Value = 46055
    # Row 13
    $ws.Range("C13").Value = 493
    $ws.Range("F13").Value = 45992
    $ws.Range("G13").Value = 30865
    $ws.Range("H13").Value = 46055
    # Row 14
    $ws.Range("C14").Value = 448
    $ws.Range("D14").Value = 422
    $ws.Range("F14").Value = 45962
    $ws.Range("H14").Value = 46055
    # Row 15
    $ws.Range("C15").Value = 409
    $ws.Range("F15").Value = 45962
    $ws.Range("G15").Value = 30834
    $ws.Range("H15").Value = 46055
    # Row 16
    $ws.Range("D16").Value = 436
    $ws.Range("H16").Value = 46055
    # Row 17
    $ws.Range("C17").Value = 396
    $ws.Range("D17").Value = 420
    $ws.Range("F17").Value = 45992
    $ws.Range("H17").Value = 46055
    # Row 18
    $ws.Range("D18").Value = 286
    $ws.Range("E18").Value = 30803
    $ws.Range("F18").Value = 45992
    $ws.Range("H18").Value = 46055
    # Row 19
    $ws.Range("D19").Value = 424
    $ws.Range("E19").Value = 30803
    $ws.Range("F19").Value = 45992
    $ws.Range("H19").Value = 46055
    # Row 20
    $ws.Range("E20").Value = 28460
    $ws.Range("F20").Value = 45992
    $ws.Range("G20").Value = 30865
    $ws.Range("H20").Value = 46055
    # Row 21
    $ws.Range("C21").Value = 349
    $ws.Range("D21").Value = 393
    $ws.Range("F21").Value = 45992
    $ws.Range("H21").Value = 46055
    # Row 22
    $ws.Range("E22").Value = 30803
    $ws.Range("F22").Value = 45992
    $ws.Range("G22").Value = 30865
    $ws.Range("H22").Value = 46055
    # Row 23
    $ws.Range("D23").Value = 406
    $ws.Range("E23").Value = 30773
    $ws.Range("F23").Value = 45962
    $ws.Range("H23").Value = 46055
    # Row 24
    $ws.Range("D24").Value = 243
    $ws.Range("H24").Value = 46055
    # Row 25
    $ws.Range("C25").Value = 420
    $ws.Range("F25").Value = 45992
    $ws.Range("G25").Value = 30865
    $ws.Range("H25").Value = 46055
    # Row 26
    $ws.Range("E26").Value = 30803
    $ws.Range("F26").Value = 45992
    $ws.Range("G26").Value = 30865
    $ws.Range("H26").Value = 46055
    # Row 28
    $ws.Range("D28").Value = 436
    $ws.Range("E28").Value = 30803
    $ws.Range("F28").Value = 45992
    $ws.Range("H28").Value = 46055
    # Row 29
    $ws.Range("C29").Value = 385
    $ws.Range("D29").Value = 393
    $ws.Range("F29").Value = 45992
    $ws.Range("H29").Value = 46055

Write-Host "Applied Global_M2 Long28_DataComp refresh."
